$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Exception Tweaking and Validation"
# The Sub_Component header in B1 had a typo/rename applied.
$ws.Range("B1").Value = "Sub_Componentt2"
